$d = $word.ActiveDocument

function Find-ParagraphIndexByText($doc, $targetText) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx += 1
        $t = $p.Range.Text
        if ($t.StartsWith($targetText)) {
            return $idx
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) "Dagkoers aanvullen" heading paragraph: split the heading text into
#    three runs ("Dagkoers" / " " / "aanvullen") plus a line break, move the
#    picture + bookmark into their own (new) paragraph, drop the stray
#    <w:lastRenderedPageBreak/>, tag the drawing with wp14 anchor/edit ids,
#    and put the page break after the picture instead of before it.
# ---------------------------------------------------------------------------

$headingIdx = Find-ParagraphIndexByText $d "Dagkoers aanvullen"
$headingPara = $d.Paragraphs.Item($headingIdx)
$headingRange = $headingPara.Range

$ns = ' xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"' +
      ' xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing"' +
      ' xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing"' +
      ' xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"' +
      ' xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"' +
      ' xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"'

$para1 = "<w:p$ns><w:pPr><w:jc w:val=`"center`"/><w:rPr><w:rStyle w:val=`"Heading3Char`"/></w:rPr></w:pPr>" +
         "<w:r><w:t>Dagkoers</w:t></w:r>" +
         "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" +
         "<w:r><w:t>aanvullen</w:t></w:r>" +
         "<w:r><w:br/></w:r></w:p>"

$drawing = '<wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="1DBE502C" wp14:editId="291665E0">' +
           '<wp:extent cx="4836461" cy="5822950"/>' +
           '<wp:effectExtent l="19050" t="19050" r="21590" b="25400"/>' +
           '<wp:docPr id="6" name="Picture 6"/>' +
           '<wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr>' +
           '<a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">' +
           '<a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
           '<pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture">' +
           '<pic:nvPicPr><pic:cNvPr id="0" name="Picture 1"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr>' +
           '<pic:blipFill><a:blip r:embed="rId7"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill>' +
           '<pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="4841122" cy="5828562"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/>' +
           '<a:ln w="6350"><a:solidFill><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:solidFill></a:ln></pic:spPr>' +
           '</pic:pic></a:graphicData></a:graphic></wp:inline>'

$para2 = "<w:p$ns><w:pPr><w:jc w:val=`"center`"/></w:pPr>" +
         "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/>" +
         "<w:r><w:rPr><w:noProof/></w:rPr><w:drawing>$drawing</w:drawing></w:r>" +
         "<w:bookmarkEnd w:id=`"0`"/>" +
         "<w:r><w:br w:type=`"page`"/></w:r></w:p>"

$headingRange.InsertXML($para1 + $para2)

# Re-locate the (now split) heading paragraph and apply the Heading 3
# character style run-by-run, since InsertXML silently drops an <w:rStyle>
# that lives directly under a run's <w:rPr> (it is kept fine under a
# paragraph mark's <w:pPr>/<w:rPr>, which is why that one was inlined above).
$headingIdx = Find-ParagraphIndexByText $d "Dagkoers"
$headingPara = $d.Paragraphs.Item($headingIdx)
$base = $headingPara.Range.Start

$d.Range($base + 0, $base + 8).Style  = "Heading 3 Char"   # "Dagkoers"
$d.Range($base + 8, $base + 9).Style  = "Heading 3 Char"   # " "
$d.Range($base + 9, $base + 18).Style = "Heading 3 Char"   # "aanvullen"
$d.Range($base + 18, $base + 19).Style = "Heading 3 Char"  # line break

# The picture paragraph that now follows holds the trailing page break;
# give that run the Heading 3 character style too.
$picturePara = $d.Paragraphs.Item($headingIdx + 1)
$pbase = $picturePara.Range.Start
$d.Range($pbase + 0, $pbase + 1).Style = "Heading 3 Char"  # page break

# ---------------------------------------------------------------------------
# 2) Drop the stray <w:lastRenderedPageBreak/> in front of "Aan- en verkopen
#    op basis van 1 indicator" (no visible text change).
# ---------------------------------------------------------------------------

$indicatorIdx = Find-ParagraphIndexByText $d "Aan- en verkopen op basis van 1 indicator"
$indicatorPara = $d.Paragraphs.Item($indicatorIdx)
$indicatorRange = $indicatorPara.Range

$indicatorXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
                '<w:r><w:t>Aan- en verkopen op basis van 1 indicator</w:t></w:r></w:p>'

$indicatorRange.InsertXML($indicatorXml)

Write-Output "edit applied"
